$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update individual cells that flip between the "NaN" shared text label
# (string index 68, text "NaN") and a literal numeric value of 1. ---

# Cells that were the text "NaN" and become the number 1
$ws.Range("AD18").Value = 1
$ws.Range("BQ19").Value = 1
$ws.Range("CJ36").Value = 1
$ws.Range("DK53").Value = 1

# Cells that were numbers and become the text "NaN"
$ws.Range("CJ159").Value = "NaN"
$ws.Range("AI172").Value = "NaN"
$ws.Range("AI173").Value = "NaN"
$ws.Range("AD186").Value = "NaN"

# --- Append a new data row (row 197) with the next day's Colombia case data ---

$row197 = @(44091,743945,2728,100411,66296,247233,27800,5552,4380,7387,7866,16649,3900,22452,29648,7057,8509,14259,12545,16531,14112,3502,2328,8872,26247,13417,10245,55696,1690,697,642,464,500,339,557,2016,4689,37392,8674,2531,43392,1058,22282,1509,9675,1629,1597,6934,1829,954,2487,2657,58595,13496,5080,8929,5924,277,1443,2667,741,2133,9388,9296,10049,14176,1939,891,12286,9878,11654,2332,1943,4914,4341,1693,5507,3190,1870,892,2732,2168,1743,1407,5988,1926,1385,1646,1995,1988,2359,1485,1202,1175,866,3300,1349,901,955,1650,1498,724,848,1218,1475,1368,1411,1115,334,362,780,730,466,536,372,657,744,525,490,372,518,132061,315453,15853,136270,83993,40408,11389)

$newRow = 197
for ($c = 1; $c -le $row197.Length; $c++) {
    $ws.Cells.Item($newRow, $c).Value = $row197[$c - 1]
}

# --- Scroll the frozen bottom-right pane down to the newly added data ---
$ws.Range("DG172").Select()
